$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 6-9 (Q4..Q7 category rows), shrinking table from A1:G9 to A1:G5
$ws.Range("A6:G9").Delete()

# Update the values for the remaining data rows (rows 2-5)
$ws.Range("B2").Value = 0.4524751023004333
$ws.Range("C2").Value = 0.7399800961295603
$ws.Range("D2").Value = 0.959507423691378
$ws.Range("E2").Value = 0.9795444980660031
$ws.Range("F2").Value = 0.901572637837064
$ws.Range("G2").Value = 14

$ws.Range("B3").Value = 0.06200830852274453
$ws.Range("C3").Value = 0.6054022942997469
$ws.Range("D3").Value = 0.6570257936368813
$ws.Range("E3").Value = 0.810571276099074
$ws.Range("F3").Value = 0.8519133780372987
$ws.Range("G3").Value = 10

$ws.Range("B4").Value = -0.06427704427340604
$ws.Range("C4").Value = 0.2800168996375158
$ws.Range("D4").Value = 0.1697350109970526
$ws.Range("E4").Value = 0.411989090871412
$ws.Range("F4").Value = 0.445784888810548
$ws.Range("G4").Value = 6

$ws.Range("B5").Value = 0.2218614552644835
$ws.Range("C5").Value = 0.2564559007953193
$ws.Range("D5").Value = 0.1149921343848131
$ws.Range("E5").Value = 0.3391049017410587
$ws.Range("F5").Value = 0.3626834130553496
$ws.Range("G5").Value = 2
